# "improved ui for first task"
# Insert a new column G (a relocated copy of the prior column F), replace
# column F with new content, make a couple of small text corrections
# elsewhere, and extend the title merge / column widths to cover the new
# column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture existing column F values (before they get overwritten) so we
#     can relocate them into the new column G. Value2 (unlike Value, which
#     misbehaves through this bridge) reliably returns the real cell text,
#     including embedded line breaks. ---
$old_F2 = $ws.Range("F2").Value2
$old_F3 = $ws.Range("F3").Value2
$old_F4 = $ws.Range("F4").Value2
$old_F5 = $ws.Range("F5").Value2
$old_F6 = $ws.Range("F6").Value2
$old_F7 = $ws.Range("F7").Value2

# --- Give column G the same width as the other data columns (30 chars).
#     ColumnWidth is expressed a bit differently than the stored column
#     width, so use the equivalent value that round-trips to 30. ---
$ws.Columns("G").ColumnWidth = 29.17

# --- Copy column F's formatting into column G so the new column matches
#     the look of the rest of the table. ---
$ws.Range("F2").Copy()
$ws.Range("G2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("F3:F8").Copy()
$ws.Range("G3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0

# --- Populate new column G with the relocated values (G8 gets slightly
#     reworded text to match the new content) ---
$ws.Range("G2").Value = $old_F2
$ws.Range("G3").Value = $old_F3
$ws.Range("G4").Value = $old_F4
$ws.Range("G5").Value = $old_F5
$ws.Range("G6").Value = $old_F6
$ws.Range("G7").Value = $old_F7
$ws.Range("G8").Value = "Experion PKS Servers with UDC Controller Collection checked and C200/C300 controllers"

# --- Replace column F with its new content ---
$ws.Range("F2").Value = "BW2024-08"
$ws.Range("F3").Value = "Experion PKS"
$ws.Range("F4").Value = "Unit Operations Controller (UOC)`nEthernet Interface Module (EIM)`nELCN Bridge`nELCN Node"
$ws.Range("F5").Value = "R511.5 initial release to R511.5 TCU5 HF2`nR520.2 initial release to R520.2 TCU6 HF3`nR530 initial release"
$ws.Range("F6").Value = "Experion PKS R511.5 TCU6 (Q4, 2024)`nExperion PKS R520.2 TCU7 (Released)`nExperion PKS R530 TCU1 (Q4, 2024)"
$ws.Range("F7").Value = "1-G9ENCXT"
$ws.Range("F8").Value = "Experion PKS customers using UOC, EIM, ELCN BRIDGE, or ELCN NODE"

# --- Minor text corrections elsewhere on the sheet ---
$ws.Range("C5").Value = "Firmware versions released before PAR 1-FM396Q5/REUCN-7883 are fixed"
$ws.Range("C6").Value = "Experion PKS R511.5 TCU6`nExperion PKS R520.2 TCU6 HF2`nExperion PKS R520.2 TCU7"

# --- Extend the title merge to cover the new column ---
$ws.Range("A1:F1").UnMerge()
$ws.Range("A1:G1").Merge()
# Merging redistributes the header's box border across the merged range;
# restore the original single-cell appearance (border only on A1).
$ws.Range("B1:G1").ClearFormats()
$ws.Range("A1").Borders.LineStyle = 1
